$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new data block in rows 12-15 (mirrors the header/series layout of rows 1-4) ---

# Row 12: header-like sequence 0..16 in columns B..R (no value in column A, like row 1)
$row12 = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16)
for ($i = 0; $i -lt $row12.Length; $i++) {
    $ws.Cells.Item(12, $i + 2).Value = $row12[$i]
}

# Row 13
$ws.Cells.Item(13, 1).Value = 1
$row13 = @(179.87795299999999,183.35587100000001,176.681567,177.056175,181.25053800000001,178.07139799999999,179.26929799999999,185.13461599999999,168.26872499999999,169.01170400000001,172.56851,181.375461,176.92695599999999,176.48434399999999,176.464879,182.13941,178.56106600000001)
for ($i = 0; $i -lt $row13.Length; $i++) {
    $ws.Cells.Item(13, $i + 2).Value = $row13[$i]
}

# Row 14
$ws.Cells.Item(14, 1).Value = 2
$row14 = @(34.884045,34.712730000000001,40.416818999999997,46.175725999999997,48.567855000000002,44.938367,46.302512999999998,71.035340000000005,64.556528999999998,65.428647999999995,66.771720999999999,63.314881999999997,68.232256000000007,69.510092999999998,69.109575000000007,73.742369999999994,69.445167999999995)
for ($i = 0; $i -lt $row14.Length; $i++) {
    $ws.Cells.Item(14, $i + 2).Value = $row14[$i]
}

# Row 15 (partial row - only A and B filled, as in the source edit)
$ws.Cells.Item(15, 1).Value = 4
$ws.Cells.Item(15, 2).Value = 71.747335000000007

# --- Move / resize the existing chart (same chart, same formatting - just repositioned) ---
$co = $ws.ChartObjects(1)
$co.Left = 892.6875
$co.Top = 297.5
$co.Width = 596.875
$co.Height = 297.5

# --- Update the active selection shown in the sheet view ---
$ws.Range("L23").Select()
